$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.717.21"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.599.48"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'211.40"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'0.510"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "'19.52"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.823.87"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "1.605.67"
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "'65.30"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "26.690.75"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").Value = "'7.23"
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "'208.82"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'2.31"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'8.92"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'142.34"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'7.10"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").Value = "1.291.41"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("E35").Value = "  -5.38%  "
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  +20.52%  "
$ws.Range("D40").Value = "'0.825"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'0.783"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "'63.10"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").Value = "1.735.99"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'91.31"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").Value = "'1.57"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'7.40"
$ws.Range("E51").Value = "  -0.80%  "
